# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" for both
#   locales (zh-cn, de-de), on the Overview sheet and on each locale's
#   own status sheet.
# - The handoff timestamps advance a bit (new handoff just generated).
# - The Status column is wider now to fit the new "Ready for handoff" text,
#   on the Overview sheet (zh-cn / de-de columns) and on the per-locale
#   sheets (Status column).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---------------
$overview.Range("E2").Value = "Ready for handoff"   # zh-cn status
$overview.Range("F2").Value = "Ready for handoff"   # de-de status
$zhcn.Range("C2").Value = "Ready for handoff"        # zh-cn sheet Status
$dede.Range("C2").Value = "Ready for handoff"        # de-de sheet Status

# --- Timestamps -----------------------------------------------------------
# Overview "Latest HO Xliff Generate Date" and de-de's "Latest Handoff
# Datetime" both advance to 10:42:40.
$overview.Range("G2").Value = "2016-09-03 10:42:40"
$dede.Range("H2").Value = "2016-09-03 10:42:40"

# zh-cn's "Latest Handoff Datetime" advances to 10:42:35.
$zhcn.Range("H2").Value = "2016-09-03 10:42:35"

# --- Column widths ----------------------------------------------------
# The Status-bearing columns grow to fit "Ready for handoff". The COM
# ColumnWidth setter snaps to whole-pixel increments, so use an input
# value that lands on the closest pixel boundary to the wider column.
$overview.Columns.Item(5).ColumnWidth = 16.33   # zh-cn column (Overview)
$overview.Columns.Item(6).ColumnWidth = 16.33   # de-de column (Overview)
$zhcn.Columns.Item(3).ColumnWidth = 16.33        # Status column (zh-cn)
$dede.Columns.Item(3).ColumnWidth = 16.33        # Status column (de-de)
